# Update the dSF column (F) values for specific rows to reflect the
# repulled data / recalculated mean, as described in the commit message
# "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -1
    4  = 3
    6  = 1
    8  = -1
    10 = 1
    11 = 2
    15 = 1
    17 = 0
    21 = 3
    22 = -2
    23 = 4
    24 = -1
    25 = -1
    27 = 1
    29 = -7
    30 = -2
    32 = -1
    34 = 1
    35 = -2
    37 = 0
    38 = -2
    39 = -2
    40 = 1
    44 = -1
    46 = -3
    47 = -10
    48 = 1
    49 = -8
    51 = -1
    58 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
